# Arrumando IA do pitstop
# Update the classification sheet with new driver data and add 4 new rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Data: Row, Name, Nationality, Team, Number, Tyre_Compound, Lap_Time
$data = @(
    @(2,  "Esteban Ocon",        "France",       "Renault",      31, "Soft", "0:01:20.271098"),
    @(3,  "Lando Norris",        "England",      "McLaren",       4, "Soft", "0:01:20.625276"),
    @(4,  "Mick Schumacher",     "Germany",      "Haas",         47, "Soft", "0:01:20.642260"),
    @(5,  "Carlos Sainz Jr",     "Spain",        "Ferrari",      55, "Soft", "0:01:20.657936"),
    @(6,  "Lewis Hamilton",      "England",      "Mercedes",     44, "Soft", "0:01:20.744348"),
    @(7,  "Valteri Bottas",      "Finland",      "Mercedes",     77, "Soft", "0:01:20.769499"),
    @(8,  "Fernando Alonso",     "Spain",        "Renault",      14, "Soft", "0:01:20.791443"),
    @(9,  "Sebastian Vettel",    "Germany",      "Aston Martin",  5, "Soft", "0:01:20.935625"),
    @(10, "Charles Lecerc",      "Monaco",       "Ferrari",      16, "Soft", "0:01:20.942232"),
    @(11, "Daniel Ricciardo",    "Australia",    "McLaren",       3, "Soft", "0:01:21.024493"),
    @(12, "Nikita Mazepin",      "Neutral",      "Haas",          9, "Soft", "0:01:21.147089"),
    @(13, "Pierre Gasly",        "France",       "Alpha Tauri",  10, "Soft", "0:01:21.238110"),
    @(14, "Sergio Pérez",        "Mexico",       "Red Bull",     11, "Soft", "0:01:21.392578"),
    @(15, "Lance Stroll",        "Canada",       "Aston Martin", 18, "Soft", "0:01:21.507750"),
    @(16, "George Russel",       "England",      "Williams",     63, "Soft", "0:01:21.580007"),
    @(17, "Kimi Raikkonen",      "Finland",      "Alfa-Romeo",    7, "Soft", "0:01:21.585747"),
    @(18, "Yuki Tsunoda",        "Japan",        "Alpha Tauri",  22, "Soft", "0:01:21.685020"),
    @(19, "Nicholas Latifi",     "Canada",       "Williams",      6, "Soft", "0:01:21.782664"),
    @(20, "Max Verstappen",      "Netherlands",  "Red Bull",     33, "Soft", "0:01:21.842005"),
    @(21, "Antonio Giovinazzi",  "Italy",        "Alfa-Romeo",   99, "Soft", "0:01:21.861092")
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 1).Value = $row[1]
    $ws.Cells.Item($r, 2).Value = $row[2]
    $ws.Cells.Item($r, 3).Value = $row[3]
    $ws.Cells.Item($r, 4).Value = $row[4]
    $ws.Cells.Item($r, 5).Value = $row[5]
    $ws.Cells.Item($r, 6).Value = $row[6]
}
